$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''257.38'
$ws.Range("E2").Value = '''4.86%'
$ws.Range("D4").Value = '''5.209'
$ws.Range("E4").Value = '''-1.17%'
$ws.Range("D5").Value = '''0.05915'
$ws.Range("E5").Value = '''3.61%'
$ws.Range("D6").Value = '''6.671'
$ws.Range("E6").Value = '''0.48%'
$ws.Range("D7").Value = '''0.8655'
$ws.Range("E7").Value = '''1.74%'
$ws.Range("D8").Value = '''1.014'
$ws.Range("E8").Value = '''14.91%'
$ws.Range("D9").Value = '''0.1417'
$ws.Range("E9").Value = '''2.49%'
$ws.Range("B10").Value = '''MandalaExchangeToken'
$ws.Range("C10").Value = '''https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '''0.07183'
$ws.Range("E10").Value = '''1.37%'
$ws.Range("B11").Value = '''BitrueCoin'
$ws.Range("C11").Value = '''https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = '''0.03164'
$ws.Range("E11").Value = '''0.74%'
$ws.Range("B12").Value = '''BitMartToken'
$ws.Range("C12").Value = '''https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Value = '''0.09224'
$ws.Range("E12").Value = '''0.03%'
$ws.Range("B13").Value = '''BitForexToken'
$ws.Range("C13").Value = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Value = '''0.001541'
$ws.Range("E13").Value = '''1.03%'
$ws.Range("B14").Value = '''One'
$ws.Range("C14").Value = '''https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '''0.0006061'
$ws.Range("E14").Value = '''1.74%'
$ws.Range("B15").Value = '''TigerCash'
$ws.Range("C15").Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.005912'
$ws.Range("E15").Value = '''-2.48%'
$ws.Range("B16").Value = '''LEO'
$ws.Range("C16").Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.478'
$ws.Range("E16").Value = '''-0.43%'
$ws.Range("B17").Value = '''GateToken'
$ws.Range("C17").Value = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.266'
$ws.Range("E17").Value = '''1.80%'
$ws.Range("B18").Value = '''BTSEToken'
$ws.Range("C18").Value = '''https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.224'
$ws.Range("E18").Value = '''1.63%'
$ws.Range("B19").Value = '''BitpandaEcosystemToken'
$ws.Range("C19").Value = '''https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3150'
$ws.Range("E19").Value = '''-0.55%'
$ws.Range("B20").Value = '''LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = '''https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '''0.03598'
$ws.Range("E20").Value = '''8.27%'
$ws.Range("E21").Value = '''0.06%'
$ws.Range("D22").Value = '''3.520'
$ws.Range("E22").Value = '''0.02%'
$ws.Range("D23").Value = '''0.04161'
$ws.Range("E23").Value = '''2.03%'
$ws.Range("D24").Value = '''0.1399'
$ws.Range("E24").Value = '''1.50%'
$ws.Range("E25").Value = '''-0.56%'
$ws.Range("E26").Value = '''8.73%'
$ws.Range("E27").Value = '''0.01%'
$ws.Range("D28").Value = '''0.0001938'
$ws.Range("E28").Value = '''33.81%'
$ws.Range("D40").Value = '''0.03827'
$ws.Range("E40").Value = '''1.01%'
$ws.Range("B41").Value = '''KickToken'
$ws.Range("C41").Value = '''https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.005508'
$ws.Range("E41").Value = '''6.66%'
$ws.Range("B42").Value = '''BKEXToken'
$ws.Range("C42").Value = '''https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1101'
$ws.Range("E42").Value = '''3.22%'
$ws.Range("D43").Value = '''0.001900'
$ws.Range("E43").Value = '''-13.63%'
$ws.Range("D44").Value = '''0.01078'
$ws.Range("E44").Value = '''13.68%'
$ws.Range("D45").Value = '''0.00005432'
$ws.Range("E45").Value = '''2.89%'
$ws.Range("E46").Value = '''0.01%'
$ws.Range("E47").Value = '''3.95%'
$ws.Range("D48").Value = '''0.002175'
$ws.Range("E48").Value = '''-4.13%'
$ws.Range("D49").Value = '''0.00002099'
$ws.Range("E49").Value = '''0.01%'
$ws.Range("D50").Value = '''0.0002000'
$ws.Range("E50").Value = '''0.01%'
